# LS8000_template.docx edit script
# - Merge the "$ 0.90 PER FT" runs into a single run
# - Replace the two tabs right after "Delivery: " with a {{lead_time}} run
# - Clear the fill on the logo pictures in the primary / first-page headers
# - Lower-case several built-in style display names

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "    $ 0." + "9" + "0 PER FT"  ->  single run "    $ 0.90 PER FT"
#    The visible text is unchanged; replacing it with itself makes the
#    engine rewrite the paragraph's run as one contiguous run.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("`$ 0.90 PER FT", $true, $false, $false, $false, $false, $true, 1, $false, "`$ 0.90 PER FT", 2)

# ---------------------------------------------------------------------------
# 2) "Delivery: " <tab><tab><tab><tab><tab><tab><tab> "FCA, Factory; Houston, TX "
#    ->
#    "Delivery: " "{{lead_time}}" <tab>x5 "FCA, Factory; Houston, TX "
# ---------------------------------------------------------------------------

# 2a. Drop two of the seven tabs that separate the label from the address.
$null = $d.Content.Find.Execute("`t`t`t`t`t`t`tFCA", $true, $false, $false, $false, $false, $true, 1, $false, "`t`t`t`t`tFCA", 2)

# 2b. Insert the {{lead_time}} placeholder right after "Delivery: ".
$rngDelivery = $d.Content
$null = $rngDelivery.Find.Execute("Delivery: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngDelivery.Collapse(0)
$rngDelivery.InsertAfter("{{lead_time}}")

# 2c. Re-select "Delivery: " and nudge formatting so it becomes its own run
#     (leaves the character formatting unchanged since Bold is toggled back).
$rngLabel = $d.Content
$null = $rngLabel.Find.Execute("Delivery: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngLabel.Font.Bold = 1
$rngLabel.Font.Bold = 0

# 2d. Same nudge for "{{lead_time}}" so it also becomes its own run, separate
#     from the trailing tabs / address text.
$rngLead = $d.Content
$null = $rngLead.Find.Execute("{{lead_time}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngLead.Font.Bold = 1
$rngLead.Font.Bold = 0

# ---------------------------------------------------------------------------
# 3) Clear the fill on the company-logo pictures in the primary and
#    first-page headers (adds <a:noFill/> to the picture shape properties).
# ---------------------------------------------------------------------------
$section = $d.Sections.Item(1)

$primaryHeader = $section.Headers.Item(1)   # wdHeaderFooterPrimary   -> header2.xml
if ($primaryHeader.Range.InlineShapes.Count -gt 0) {
    $primaryHeader.Range.InlineShapes.Item(1).Fill.Visible = $false
}

$firstPageHeader = $section.Headers.Item(2) # wdHeaderFooterFirstPage -> header3.xml
if ($firstPageHeader.Range.InlineShapes.Count -gt 0) {
    $firstPageHeader.Range.InlineShapes.Item(1).Fill.Visible = $false
}

# ---------------------------------------------------------------------------
# 4) Lower-case the display name of several built-in styles.
# ---------------------------------------------------------------------------
$d.Styles.Item("Heading1").NameLocal = "heading 1"
$d.Styles.Item("Heading2").NameLocal = "heading 2"
$d.Styles.Item("Caption").NameLocal = "caption"
$d.Styles.Item("Header").NameLocal = "header"
$d.Styles.Item("Footer").NameLocal = "footer"

Write-Host "Edit complete"
